$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18,8).Value = 2922.111
$ws.Cells.Item(18,9).Value = 2412.375
$ws.Cells.Item(18,10).Value = 7000
$ws.Cells.Item(18,11).Value = 2412.375
$ws.Cells.Item(18,12).Value = 7000
$ws.Cells.Item(18,13).Value = -2128.375
$ws.Cells.Item(18,14).Value = -7568

$ws.Cells.Item(40,8).Value = 3614.1428
$ws.Cells.Item(40,9).Value = 14999
$ws.Cells.Item(40,10).Value = 1716.6666
$ws.Cells.Item(40,11).Value = 14999
$ws.Cells.Item(40,12).Value = 1716.6666
$ws.Cells.Item(40,13).Value = -14824
$ws.Cells.Item(40,14).Value = -2066.6666

$ws.Cells.Item(43,8).Value = 14475.444
$ws.Cells.Item(43,9).Value = 5399.5
$ws.Cells.Item(43,10).Value = 21736.2
$ws.Cells.Item(43,11).Value = 5399.5
$ws.Cells.Item(43,12).Value = 21736.2
$ws.Cells.Item(43,13).Value = -5330.5
$ws.Cells.Item(43,14).Value = -21874.2

$ws.Cells.Item(51,8).Value = 3799
$ws.Cells.Item(51,9).Value = 0
$ws.Cells.Item(51,10).Value = 3799
$ws.Cells.Item(51,11).Value = 0
$ws.Cells.Item(51,12).Value = 3799
$ws.Cells.Item(51,13).ClearContents()
$ws.Cells.Item(51,14).Value = -4767

$ws.Cells.Item(55,8).Value = 243.63158
$ws.Cells.Item(55,9).Value = 299
$ws.Cells.Item(55,10).Value = 203.36363
$ws.Cells.Item(55,11).Value = 299
$ws.Cells.Item(55,12).Value = 203.36363
$ws.Cells.Item(55,13).Value = -85
$ws.Cells.Item(55,14).Value = -631.3636300000001

$ws.Cells.Item(64,8).Value = 3012.9583
$ws.Cells.Item(64,10).Value = 2880.6667
$ws.Cells.Item(64,12).Value = 2880.6667
$ws.Cells.Item(64,14).Value = -3376.6667

$ws.Cells.Item(67,8).Value = 3012.9583
$ws.Cells.Item(67,10).Value = 2880.6667
$ws.Cells.Item(67,12).Value = 2880.6667
$ws.Cells.Item(67,14).Value = -4596.6667

$ws.Cells.Item(74,8).Value = 3324.2632
$ws.Cells.Item(74,9).Value = 3166.111
$ws.Cells.Item(74,11).Value = 3166.111
$ws.Cells.Item(74,13).Value = -2230.111

$ws.Cells.Item(77,8).Value = 3324.2632
$ws.Cells.Item(77,9).Value = 3166.111
$ws.Cells.Item(77,11).Value = 15830.555
$ws.Cells.Item(77,13).Value = -11150.555

$ws.Cells.Item(111,8).Value = 1770.0952
$ws.Cells.Item(111,9).Value = 2651.3333
$ws.Cells.Item(111,10).Value = 1417.6
$ws.Cells.Item(111,11).Value = 7953.999899999999
$ws.Cells.Item(111,12).Value = 4252.799999999999
$ws.Cells.Item(111,13).Value = -4886.999899999999
$ws.Cells.Item(111,14).Value = -10386.8

$ws.Cells.Item(137,8).Value = 5408311.5
$ws.Cells.Item(137,9).Value = 1811.2222
$ws.Cells.Item(137,11).Value = 5433.6666
$ws.Cells.Item(137,13).Value = -2883.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88,8).Value = 2433.7334
$ws.Cells.Item(88,9).Value = 2084.3333
$ws.Cells.Item(88,11).Value = 2084.3333
$ws.Cells.Item(88,13).Value = -1678.3333

$ws.Cells.Item(91,8).Value = 2433.7334
$ws.Cells.Item(91,9).Value = 2084.3333
$ws.Cells.Item(91,11).Value = 2084.3333
$ws.Cells.Item(91,13).Value = -680.3332999999998

$ws.Cells.Item(124,8).Value = 27809.666
$ws.Cells.Item(124,10).Value = 27809.666
$ws.Cells.Item(124,12).Value = 27809.666
$ws.Cells.Item(124,14).Value = -37629.666

$ws.Cells.Item(125,8).Value = 48377.918
$ws.Cells.Item(125,10).Value = 48377.918
$ws.Cells.Item(125,12).Value = 48377.918
$ws.Cells.Item(125,14).Value = -58217.918

$ws.Cells.Item(132,8).Value = 6252176.5
$ws.Cells.Item(132,9).Value = 10871272
$ws.Cells.Item(132,10).Value = 2813
$ws.Cells.Item(132,11).Value = 32613816
$ws.Cells.Item(132,12).Value = 8439
$ws.Cells.Item(132,13).Value = -32611286
$ws.Cells.Item(132,14).Value = -13499

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107,8).Value = 1548.4
$ws.Cells.Item(107,9).Value = 1548.4
$ws.Cells.Item(107,10).Value = 0
$ws.Cells.Item(107,11).Value = 1548.4
$ws.Cells.Item(107,12).Value = 0
$ws.Cells.Item(107,13).Value = 371.5999999999999
$ws.Cells.Item(107,14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(47,8).Value = 40071
$ws.Cells.Item(47,9).Value = 0
$ws.Cells.Item(47,11).Value = 0
$ws.Cells.Item(47,13).ClearContents()

$ws.Cells.Item(107,8).Value = 308.0465
$ws.Cells.Item(107,9).Value = 285.44
$ws.Cells.Item(107,10).Value = 339.44446
$ws.Cells.Item(107,11).Value = 285.44
$ws.Cells.Item(107,12).Value = 339.44446
$ws.Cells.Item(107,13).Value = 1634.56
$ws.Cells.Item(107,14).Value = -4179.44446

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68,8).Value = 902.1875
$ws.Cells.Item(68,9).Value = 615.5263
$ws.Cells.Item(68,10).Value = 1612.6086
$ws.Cells.Item(68,11).Value = 1846.5789
$ws.Cells.Item(68,12).Value = 4837.825800000001
$ws.Cells.Item(68,13).Value = -1035.5789
$ws.Cells.Item(68,14).Value = -6459.825800000001

$ws.Cells.Item(71,8).Value = 902.1875
$ws.Cells.Item(71,9).Value = 615.5263
$ws.Cells.Item(71,10).Value = 1612.6086
$ws.Cells.Item(71,11).Value = 5539.736699999999
$ws.Cells.Item(71,12).Value = 14513.4774
$ws.Cells.Item(71,13).Value = -1483.736699999999
$ws.Cells.Item(71,14).Value = -22625.4774

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(111,8).Value = 25146
$ws.Cells.Item(111,10).Value = 25146
$ws.Cells.Item(111,12).Value = 25146
$ws.Cells.Item(111,14).Value = -31280

$ws.Cells.Item(117,8).Value = 59310
$ws.Cells.Item(117,10).Value = 59310
$ws.Cells.Item(117,12).Value = 59310
$ws.Cells.Item(117,14).Value = -66194

$ws.Cells.Item(132,8).Value = 5835.7617
$ws.Cells.Item(132,9).Value = 4669.5557
$ws.Cells.Item(132,10).Value = 6710.4165
$ws.Cells.Item(132,11).Value = 14008.6671
$ws.Cells.Item(132,12).Value = 20131.2495
$ws.Cells.Item(132,13).Value = -11478.6671
$ws.Cells.Item(132,14).Value = -25191.2495

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16,8).Value = 1873
$ws.Cells.Item(16,9).Value = 1848.7778
$ws.Cells.Item(16,10).Value = 2200
$ws.Cells.Item(16,11).Value = 1848.7778
$ws.Cells.Item(16,12).Value = 2200
$ws.Cells.Item(16,13).Value = -1678.7778
$ws.Cells.Item(16,14).Value = -2540

$ws.Cells.Item(46,8).Value = 954.7692
$ws.Cells.Item(46,9).Value = 770
$ws.Cells.Item(46,10).Value = 1250.4
$ws.Cells.Item(46,11).Value = 770
$ws.Cells.Item(46,12).Value = 1250.4
$ws.Cells.Item(46,13).Value = -582
$ws.Cells.Item(46,14).Value = -1626.4

$ws.Cells.Item(68,8).Value = 1255.5186
$ws.Cells.Item(68,9).Value = 970.125
$ws.Cells.Item(68,10).Value = 1375.6842
$ws.Cells.Item(68,11).Value = 970.125
$ws.Cells.Item(68,12).Value = 1375.6842
$ws.Cells.Item(68,13).Value = -221.125
$ws.Cells.Item(68,14).Value = -2873.6842

$ws.Cells.Item(71,8).Value = 1255.5186
$ws.Cells.Item(71,9).Value = 970.125
$ws.Cells.Item(71,10).Value = 1375.6842
$ws.Cells.Item(71,11).Value = 4850.625
$ws.Cells.Item(71,12).Value = 6878.420999999999
$ws.Cells.Item(71,13).Value = -1106.625
$ws.Cells.Item(71,14).Value = -14366.421

$ws.Cells.Item(116,8).Value = 52680
$ws.Cells.Item(116,10).Value = 52680
$ws.Cells.Item(116,12).Value = 52680
$ws.Cells.Item(116,14).Value = -61858

$ws.Cells.Item(132,8).Value = 10007305
$ws.Cells.Item(132,9).Value = 4697.4688
$ws.Cells.Item(132,10).Value = 27789718
$ws.Cells.Item(132,11).Value = 14092.4064
$ws.Cells.Item(132,12).Value = 83369154
$ws.Cells.Item(132,13).Value = -11562.4064
$ws.Cells.Item(132,14).Value = -83374214

$ws.Cells.Item(136,8).Value = 28850764
$ws.Cells.Item(136,9).Value = 38463190
$ws.Cells.Item(136,10).Value = 13481.538
$ws.Cells.Item(136,11).Value = 115389570
$ws.Cells.Item(136,12).Value = 40444.614
$ws.Cells.Item(136,13).Value = -115387020
$ws.Cells.Item(136,14).Value = -45544.614
